# Update Overview Decks for July 2025 (#60)
#
# The deck's tables were using the old custom table style
# {E6D9C5E6-9F9E-41A7-B0E6-88E029BFA50D}; re-point every table in the
# presentation that still uses it to the new style
# {98510C90-4868-4D3F-B63A-26FF2F9ABE26}.

$p = $ppt.ActivePresentation

$oldStyleId = "{E6D9C5E6-9F9E-41A7-B0E6-88E029BFA50D}"
$newStyleId = "{98510C90-4868-4D3F-B63A-26FF2F9ABE26}"

$updated = 0

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)

        if ($shape.HasTable) {
            $table = $shape.Table

            # Only touch tables that are still on the old style so the
            # script is idempotent / safe even if it is re-run.
            if ($table.Style.Name -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
                $updated = $updated + 1
            }
        }
    }
}

Write-Host ("Updated {0} table(s) from {1} to {2}" -f $updated, $oldStyleId, $newStyleId)
